$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H9").Value = 117.375
$ws.Range("I9").Value = 118
$ws.Range("K9").Value = 118
$ws.Range("M9").Value = 51
$ws.Range("H33").Value = 222.72728
$ws.Range("I33").Value = 241.94737
$ws.Range("J33").Value = 101
$ws.Range("K33").Value = 241.94737
$ws.Range("L33").Value = 101
$ws.Range("M33").Value = -12.94737000000001
$ws.Range("N33").Value = -559
$ws.Range("H80").Value = 791.9545000000001
$ws.Range("I80").Value = 710.6316
$ws.Range("J80").Value = 853.76
$ws.Range("K80").Value = 2131.8948
$ws.Range("L80").Value = 2561.28
$ws.Range("M80").Value = -1133.8948
$ws.Range("N80").Value = -4557.28
$ws.Range("H82").Value = 7693105
$ws.Range("I82").Value = 8334123.5
$ws.Range("K82").Value = 25002370.5
$ws.Range("M82").Value = -25001964.5
$ws.Range("H83").Value = 791.9545000000001
$ws.Range("I83").Value = 710.6316
$ws.Range("J83").Value = 853.76
$ws.Range("K83").Value = 6395.6844
$ws.Range("L83").Value = 7683.84
$ws.Range("M83").Value = -1403.6844
$ws.Range("N83").Value = -17667.84
$ws.Range("H85").Value = 7693105
$ws.Range("I85").Value = 8334123.5
$ws.Range("K85").Value = 25002370.5
$ws.Range("M85").Value = -25000966.5
$ws.Range("H86").Value = 5498.5
$ws.Range("I86").Value = 3999
$ws.Range("J86").Value = 5998.3335
$ws.Range("K86").Value = 3999
$ws.Range("L86").Value = 5998.3335
$ws.Range("M86").Value = -2876
$ws.Range("N86").Value = -8244.333500000001
$ws.Range("H88").Value = 531221.4
$ws.Range("I88").Value = 4479.8335
$ws.Range("J88").Value = 774332.9
$ws.Range("K88").Value = 4479.8335
$ws.Range("L88").Value = 774332.9
$ws.Range("M88").Value = -4073.8335
$ws.Range("N88").Value = -775144.9
$ws.Range("H89").Value = 5498.5
$ws.Range("I89").Value = 3999
$ws.Range("J89").Value = 5998.3335
$ws.Range("K89").Value = 19995
$ws.Range("L89").Value = 29991.6675
$ws.Range("M89").Value = -14379
$ws.Range("N89").Value = -41223.6675
$ws.Range("H91").Value = 531221.4
$ws.Range("I91").Value = 4479.8335
$ws.Range("J91").Value = 774332.9
$ws.Range("K91").Value = 4479.8335
$ws.Range("L91").Value = 774332.9
$ws.Range("M91").Value = -3075.8335
$ws.Range("N91").Value = -777140.9
$ws.Range("H94").Value = 5296.7144
$ws.Range("I94").Value = 5296.7144
$ws.Range("K94").Value = 5296.7144
$ws.Range("M94").Value = -4845.7144
$ws.Range("H100").Value = 2400.8333
$ws.Range("I100").Value = 2174.25
$ws.Range("K100").Value = 2174.25
$ws.Range("M100").Value = -1633.25
$ws.Range("H138").Value = 2995.3274
$ws.Range("J138").Value = 3439.3022
$ws.Range("L138").Value = 10317.9066
$ws.Range("N138").Value = -20597.9066

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H2").Value = 2179.9333
$ws.Range("J2").Value = 4332
$ws.Range("L2").Value = 4332
$ws.Range("N2").Value = -4558
$ws.Range("H88").Value = 1521.6666
$ws.Range("I88").Value = 1204.4
$ws.Range("J88").Value = 1786.0555
$ws.Range("K88").Value = 1204.4
$ws.Range("L88").Value = 1786.0555
$ws.Range("M88").Value = -798.4000000000001
$ws.Range("N88").Value = -2598.0555
$ws.Range("H91").Value = 1521.6666
$ws.Range("I91").Value = 1204.4
$ws.Range("J91").Value = 1786.0555
$ws.Range("K91").Value = 1204.4
$ws.Range("L91").Value = 1786.0555
$ws.Range("M91").Value = 199.5999999999999
$ws.Range("N91").Value = -4594.0555
$ws.Range("H102").Value = 7243
$ws.Range("I102").Value = 8200.166999999999
$ws.Range("K102").Value = 8200.166999999999
$ws.Range("M102").Value = -6578.166999999999
$ws.Range("H110").Value = 969.75
$ws.Range("I110").Value = 1129.6666
$ws.Range("K110").Value = 1129.6666
$ws.Range("M110").Value = 915.3334
$ws.Range("H116").Value = 2179.9333
$ws.Range("J116").Value = 4332
$ws.Range("L116").Value = 4332
$ws.Range("N116").Value = -8920

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H3").Value = 2179.9333
$ws.Range("J3").Value = 4332
$ws.Range("L3").Value = 4332
$ws.Range("N3").Value = -4560
$ws.Range("H86").Value = 2880
$ws.Range("I86").Value = 2926.6667
$ws.Range("J86").Value = 2775
$ws.Range("K86").Value = 2926.6667
$ws.Range("L86").Value = 2775
$ws.Range("M86").Value = -1803.6667
$ws.Range("N86").Value = -5021
$ws.Range("H89").Value = 2880
$ws.Range("I89").Value = 2926.6667
$ws.Range("J89").Value = 2775
$ws.Range("K89").Value = 14633.3335
$ws.Range("L89").Value = 13875
$ws.Range("M89").Value = -9017.333500000001
$ws.Range("N89").Value = -25107
$ws.Range("H94").Value = 2417.6667
$ws.Range("I94").Value = 2417.6667
$ws.Range("K94").Value = 2417.6667
$ws.Range("M94").Value = -1966.6667
$ws.Range("H122").Value = 99999
$ws.Range("J122").Value = 99999
$ws.Range("L122").Value = 99999
$ws.Range("M122").Value = -109799

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H31").Value = 12298
$ws.Range("J31").Value = 15124.9
$ws.Range("L31").Value = 15124.9
$ws.Range("N31").Value = -15714.9
$ws.Range("H34").Value = 12298
$ws.Range("J34").Value = 15124.9
$ws.Range("L34").Value = 15124.9
$ws.Range("N34").Value = -15528.9

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H122").Value = 1650.125
$ws.Range("I122").Value = 1406.8572
$ws.Range("J122").Value = 1990.7
$ws.Range("K122").Value = 12661.7148
$ws.Range("L122").Value = 17916.3
$ws.Range("M122").Value = -10211.7148
$ws.Range("N122").Value = -22816.3
$ws.Range("H138").Value = 6950
$ws.Range("I138").Value = 4266.6665
$ws.Range("K138").Value = 12799.9995
$ws.Range("M138").Value = -7659.999500000002

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H107").Value = 1256.8096
$ws.Range("I107").Value = 1302.25
$ws.Range("J107").Value = 1111.4
$ws.Range("K107").Value = 1302.25
$ws.Range("L107").Value = 1111.4
$ws.Range("M107").Value = 617.75
$ws.Range("N107").Value = -4951.4
$ws.Range("H113").Value = 1678.8
$ws.Range("I113").Value = 1465.3334
$ws.Range("K113").Value = 1465.3334
$ws.Range("M113").Value = 704.6666

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H68").Value = 2617.9375
$ws.Range("I68").Value = 2356.8462
$ws.Range("K68").Value = 2356.8462
$ws.Range("M68").Value = -1607.8462
$ws.Range("H71").Value = 2617.9375
$ws.Range("I71").Value = 2356.8462
$ws.Range("K71").Value = 11784.231
$ws.Range("M71").Value = -8040.231

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H81").Value = 3394
$ws.Range("I81").Value = 2180.8333
$ws.Range("J81").Value = 6305.6
$ws.Range("K81").Value = 4361.6666
$ws.Range("L81").Value = 12611.2
$ws.Range("M81").Value = -3300.6666
$ws.Range("N81").Value = -14733.2
$ws.Range("H84").Value = 3394
$ws.Range("I84").Value = 2180.8333
$ws.Range("J84").Value = 6305.6
$ws.Range("K84").Value = 21808.333
$ws.Range("L84").Value = 63056
$ws.Range("M84").Value = -16504.333
$ws.Range("N84").Value = -73664
$ws.Range("H100").Value = 1272.3636
$ws.Range("I100").Value = 713
$ws.Range("K100").Value = 1426
$ws.Range("M100").Value = -885
$ws.Range("H101").Value = 24517.166
$ws.Range("J101").Value = 24517.166
$ws.Range("L101").Value = 24517.166
$ws.Range("N101").Value = -31007.166
$ws.Range("H136").Value = 4457
$ws.Range("I136").Value = 5102.3076
$ws.Range("J136").Value = 2359.75
$ws.Range("K136").Value = 15306.9228
$ws.Range("L136").Value = 7079.25
$ws.Range("M136").Value = -12756.9228
$ws.Range("N136").Value = -12179.25
